{"js": "// The commit re-saves the manuscript and picks up the \"Compact\" paragraph\n// style gaining an explicit 10pt (20 half-points) run size\n// (<w:rPr><w:sz w:val=\"20\"/></w:rPr>) on top of its existing\n// \"before/after 36 twips\" paragraph spacing.\n//\n// (The rest of the underlying diff \u2014 regenerated w:tmpl GUIDs on the\n// legacy list definitions, the latentStyles w:count bump, and the new\n// \"Smart Link Error\" lsdException \u2014 are Word's own save-time bookkeeping\n// for internal list-template caches / the built-in style registry. They\n// are not part of the Word JS object model, so they fall out of scope for\n// a content edit like this one and are not touched here.)\n\nconst styles = context.document.styles;\nconst compact = styles.getByNameOrNullObject(\"Compact\");\ncompact.load(\"nameLocal\");\nawait context.sync();\n\nif (!compact.isNullObject) {\n  compact.font.size = 10;\n  await context.sync();\n}\n", "ps1": "# The commit re-saves the manuscript and picks up the \"Compact\" paragraph\n# style gaining an explicit 10pt (20 half-points) run size\n# (<w:rPr><w:sz w:val=\"20\"/></w:rPr>) on top of its existing\n# \"before/after 36 twips\" paragraph spacing.\n#\n# (The rest of the underlying diff -- regenerated w:tmpl GUIDs on the\n# legacy list definitions, the latentStyles w:count bump, and the new\n# \"Smart Link Error\" lsdException -- are Word's own save-time bookkeeping\n# for internal list-template caches / the built-in style registry. They\n# are not part of the Word COM object model, so they fall out of scope\n# for a content edit like this one and are not touched here.)\n\n$d = $word.ActiveDocument\n\n$compact = $d.Styles(\"Compact\")\nif ($compact) {\n    $compact.Font.Size = 10\n}\n"}
